# Auto-generated edit script applying the Louisoix_Profits value updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 3191
$ws.Range("J32").Value = 3162.6
$ws.Range("L32").Value = 3162.6
$ws.Range("N32").Value = -3814.6

# Row 55
$ws.Range("H55").Value = 285.23077
$ws.Range("I55").Value = 253
$ws.Range("K55").Value = 253
$ws.Range("M55").Value = -39

# Row 70
$ws.Range("H70").Value = 3918.2856
$ws.Range("I70").Value = 2499
$ws.Range("K70").Value = 7497
$ws.Range("M70").Value = -7227

# Row 73
$ws.Range("H73").Value = 3918.2856
$ws.Range("I73").Value = 2499
$ws.Range("K73").Value = 7497
$ws.Range("M73").Value = -6561

# Row 107
$ws.Range("H107").Value = 1338
$ws.Range("I107").Value = 894.8889
$ws.Range("J107").Value = 2002.6666
$ws.Range("K107").Value = 894.8889
$ws.Range("L107").Value = 2002.6666
$ws.Range("M107").Value = 1025.1111
$ws.Range("N107").Value = -5842.6666

# Row 111
$ws.Range("H111").Value = 2264.353
$ws.Range("I111").Value = 2205.8333
$ws.Range("K111").Value = 6617.499899999999
$ws.Range("M111").Value = -3550.499899999999

# Row 129
$ws.Range("H129").Value = 1201.875
$ws.Range("I129").Value = 1016.4286
$ws.Range("J129").Value = 2500
$ws.Range("K129").Value = 3049.2858
$ws.Range("L129").Value = 7500
$ws.Range("M129").Value = 1950.7142
$ws.Range("N129").Value = -17500

# Row 132
$ws.Range("H132").Value = 5421.75
$ws.Range("I132").Value = 4249.472
$ws.Range("K132").Value = 12748.416
$ws.Range("M132").Value = -10218.416

# Row 137
$ws.Range("H137").Value = 3718.4
$ws.Range("I137").Value = 3790.923
$ws.Range("J137").Value = 3247
$ws.Range("K137").Value = 11372.769
$ws.Range("L137").Value = 9741
$ws.Range("M137").Value = -8822.769
$ws.Range("N137").Value = -14841


$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 16573.4
$ws.Range("I32").Value = 19351.736
$ws.Range("K32").Value = 19351.736
$ws.Range("M32").Value = -19064.736

# Row 61
$ws.Range("H61").Value = 3393.6667
$ws.Range("I61").Value = 3009.8
$ws.Range("J61").Value = 3873.5
$ws.Range("K61").Value = 3009.8
$ws.Range("L61").Value = 3873.5
$ws.Range("M61").Value = -2797.8
$ws.Range("N61").Value = -4297.5

# Row 102
$ws.Range("H102").Value = 5413.857
$ws.Range("I102").Value = 4640.6
$ws.Range("K102").Value = 4640.6
$ws.Range("M102").Value = -3018.6

# Row 110
$ws.Range("H110").Value = 3738
$ws.Range("I110").Value = 3384.3333
$ws.Range("K110").Value = 3384.3333
$ws.Range("M110").Value = -1339.3333

# Row 112
$ws.Range("H112").Value = 27999.666
$ws.Range("J112").Value = 27999.666
$ws.Range("L112").Value = 27999.666
$ws.Range("N112").Value = -30953.666

# Row 132
$ws.Range("H132").Value = 29207.475
$ws.Range("I132").Value = 52546.3
$ws.Range("K132").Value = 157638.9
$ws.Range("M132").Value = -155108.9

# Row 136
$ws.Range("H136").Value = 3393.6667
$ws.Range("I136").Value = 3009.8
$ws.Range("J136").Value = 3873.5
$ws.Range("K136").Value = 9029.400000000001
$ws.Range("L136").Value = 11620.5
$ws.Range("M136").Value = -6479.400000000001
$ws.Range("N136").Value = -16720.5


$ws = $wb.Worksheets.Item("CRP")
# Row 36
$ws.Range("H36").Value = 5666.6665
$ws.Range("I36").Value = 6000
$ws.Range("J36").Value = 5000
$ws.Range("K36").Value = 6000
$ws.Range("L36").Value = 5000
$ws.Range("M36").Value = -5612
$ws.Range("N36").Value = -5776

# Row 40
$ws.Range("H40").Value = 5666.6665
$ws.Range("I40").Value = 6000
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 6000
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -5840
$ws.Range("N40").Value = -5320

# Row 58
$ws.Range("H58").Value = 115099.445
$ws.Range("I58").Value = 115099.445
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 115099.445
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -114896.445
$ws.Range("N58").ClearContents()

# Row 94
$ws.Range("H94").Value = 1723.2941
$ws.Range("I94").Value = 1674.75
$ws.Range("J94").Value = 1766.4445
$ws.Range("K94").Value = 1674.75
$ws.Range("L94").Value = 1766.4445
$ws.Range("M94").Value = -1223.75
$ws.Range("N94").Value = -2668.4445

# Row 105
$ws.Range("H105").Value = 1581.3125
$ws.Range("I105").Value = 1593.2307
$ws.Range("J105").Value = 1529.6666
$ws.Range("K105").Value = 1593.2307
$ws.Range("L105").Value = 1529.6666
$ws.Range("M105").Value = 153.7692999999999
$ws.Range("N105").Value = -5023.6666

# Row 120
$ws.Range("H120").Value = 98749
$ws.Range("J120").Value = 98749
$ws.Range("L120").Value = 98749
$ws.Range("N120").Value = -106007

# Row 132
$ws.Range("H132").Value = 1830.5454
$ws.Range("I132").Value = 1916.6
$ws.Range("J132").Value = 970
$ws.Range("K132").Value = 5749.799999999999
$ws.Range("L132").Value = 2910
$ws.Range("M132").Value = -3219.799999999999
$ws.Range("N132").Value = -7970

# Row 134
$ws.Range("H134").Value = 53220.25
$ws.Range("I134").Value = 55705.633
$ws.Range("K134").Value = 167116.899
$ws.Range("M134").Value = -164581.899

# Row 136
$ws.Range("H136").Value = 115099.445
$ws.Range("I136").Value = 115099.445
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 345298.335
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -342748.335
$ws.Range("N136").ClearContents()


$ws = $wb.Worksheets.Item("CUL")
# Row 33
$ws.Range("H33").Value = 683.3333
$ws.Range("J33").Value = 1400
$ws.Range("L33").Value = 8400
$ws.Range("N33").Value = -8966

# Row 68
$ws.Range("H68").Value = 144971.28
$ws.Range("I68").Value = 2475
$ws.Range("J68").Value = 334966.34
$ws.Range("K68").Value = 7425
$ws.Range("L68").Value = 1004899.02
$ws.Range("M68").Value = -6614
$ws.Range("N68").Value = -1006521.02

# Row 71
$ws.Range("H71").Value = 144971.28
$ws.Range("I71").Value = 2475
$ws.Range("J71").Value = 334966.34
$ws.Range("K71").Value = 22275
$ws.Range("L71").Value = 3014697.06
$ws.Range("M71").Value = -18219
$ws.Range("N71").Value = -3022809.06

# Row 76
$ws.Range("H76").Value = 14333
$ws.Range("I76").Value = 8999
$ws.Range("J76").Value = 17000
$ws.Range("K76").Value = 26997
$ws.Range("L76").Value = 51000
$ws.Range("M76").Value = -26614
$ws.Range("N76").Value = -51766

# Row 79
$ws.Range("H79").Value = 14333
$ws.Range("I79").Value = 8999
$ws.Range("J79").Value = 17000
$ws.Range("K79").Value = 26997
$ws.Range("L79").Value = 51000
$ws.Range("M79").Value = -25671
$ws.Range("N79").Value = -53652

# Row 104
$ws.Range("H104").Value = 1237.75
$ws.Range("I104").Value = 983.6667
$ws.Range("K104").Value = 2951.0001
$ws.Range("M104").Value = -330.0001000000002

# Row 117
$ws.Range("H117").Value = 3319.6924
$ws.Range("I117").Value = 989
$ws.Range("J117").Value = 3513.9167
$ws.Range("K117").Value = 2967
$ws.Range("L117").Value = 10541.7501
$ws.Range("M117").Value = 475
$ws.Range("N117").Value = -17425.7501


$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 3533.1667
$ws.Range("I102").Value = 2275.889
$ws.Range("K102").Value = 2275.889
$ws.Range("M102").Value = -653.8890000000001

# Row 126
$ws.Range("H126").Value = 4575.2856
$ws.Range("I126").Value = 4762.923
$ws.Range("K126").Value = 14288.769
$ws.Range("M126").Value = -11818.769


$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 2081.442
$ws.Range("I16").Value = 1700.5428
$ws.Range("K16").Value = 1700.5428
$ws.Range("M16").Value = -1530.5428

# Row 100
$ws.Range("H100").Value = 15000
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

# Row 127
$ws.Range("H127").Value = 49248.75
$ws.Range("J127").Value = 49248.75
$ws.Range("L127").Value = 49248.75
$ws.Range("N127").Value = -59168.75

# Row 136
$ws.Range("H136").Value = 9591.666999999999
$ws.Range("I136").Value = 9387.5
$ws.Range("K136").Value = 28162.5
$ws.Range("M136").Value = -25612.5


$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 3537.625
$ws.Range("I100").Value = 2959.6
$ws.Range("K100").Value = 5919.2
$ws.Range("M100").Value = -5378.2

# Row 107
$ws.Range("H107").Value = 4832.5
$ws.Range("I107").Value = 3999
$ws.Range("J107").Value = 5666
$ws.Range("K107").Value = 11997
$ws.Range("L107").Value = 16998
$ws.Range("M107").Value = -10077
$ws.Range("N107").Value = -20838

# Row 122
$ws.Range("H122").Value = 7214.8184
$ws.Range("I122").Value = 8070.263
$ws.Range("K122").Value = 24210.789
$ws.Range("M122").Value = -21760.789

# Row 132
$ws.Range("H132").Value = 49739.145
$ws.Range("I132").Value = 54622.74
$ws.Range("K132").Value = 163868.22
$ws.Range("M132").Value = -161338.22

# Row 136
$ws.Range("H136").Value = 3972
$ws.Range("I136").Value = 2241
$ws.Range("J136").Value = 8299.5
$ws.Range("K136").Value = 6723
$ws.Range("L136").Value = 24898.5
$ws.Range("M136").Value = -4173
$ws.Range("N136").Value = -29998.5

